$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

# Replace whole text with the new three-run text.
$tr.Text = "Play the 24 game using the following numbers"

# Select just the "24 game" substring and add the hyperlink to it.
$start = ("Play the ").Length + 1
$len = ("24 game").Length
$linkRange = $tr.Characters($start, $len)
$linkRange.ActionSettings.Item(1).Hyperlink.Address = "https://en.wikipedia.org/wiki/24_(puzzle)"
